$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the "Price" (column E) values for every menu item (rows 2-40)
# ---------------------------------------------------------------------------
$prices = New-Object 'object[,]' 39,1
$prices[0,0]  = 14000   # E2
$prices[1,0]  = 18000   # E3
$prices[2,0]  = 20000   # E4
$prices[3,0]  = 20000   # E5
$prices[4,0]  = 22000   # E6
$prices[5,0]  = 22000   # E7
$prices[6,0]  = 22000   # E8
$prices[7,0]  = 24000   # E9
$prices[8,0]  = 24000   # E10
$prices[9,0]  = 24000   # E11
$prices[10,0] = 24000   # E12
$prices[11,0] = 24000   # E13
$prices[12,0] = 24000   # E14
$prices[13,0] = 22000   # E15
$prices[14,0] = 22000   # E16
$prices[15,0] = 22000   # E17
$prices[16,0] = 22000   # E18
$prices[17,0] = 22000   # E19
$prices[18,0] = 22000   # E20
$prices[19,0] = 22000   # E21
$prices[20,0] = 12000   # E22
$prices[21,0] = 16000   # E23
$prices[22,0] = 16000   # E24
$prices[23,0] = 6000    # E25
$prices[24,0] = 15000   # E26
$prices[25,0] = 15000   # E27
$prices[26,0] = 15000   # E28
$prices[27,0] = 12000   # E29
$prices[28,0] = 14000   # E30
$prices[29,0] = 14000   # E31
$prices[30,0] = 16000   # E32
$prices[31,0] = 20000   # E33
$prices[32,0] = 20000   # E34
$prices[33,0] = 20000   # E35
$prices[34,0] = 22000   # E36
$prices[35,0] = 22000   # E37
$prices[36,0] = 20000   # E38
$prices[37,0] = 10000   # E39
$prices[38,0] = 10000   # E40
$ws.Range("E2:E40").Value = $prices

# ---------------------------------------------------------------------------
# 2. Remove the trailing empty/template rows (41-46) - the sheet now ends at
#    row 40, shrinking the used range from A1:E46 to A1:E40.
# ---------------------------------------------------------------------------
$ws.Range("A41:E46").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3. Re-align all the table text to the left (it used to be centered / default)
#    and drop the bold weight that used to be applied to the "Name" column and
#    the header row.
# ---------------------------------------------------------------------------
$usedRange = $ws.Range("A1:E40")
$usedRange.HorizontalAlignment = -4131   # xlLeft
$ws.Range("A1:E1").Font.Bold = $false
$ws.Range("A2").Font.Bold = $false
$ws.Range("A3:A40").Font.Bold = $false

# ---------------------------------------------------------------------------
# 4. Small row-height tweaks that come from the content/format change above.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 16.5
$ws.Rows.Item(2).RowHeight = 30.6
$ws.Rows.Item(3).RowHeight = 44.45
$ws.Rows.Item(4).RowHeight = 44.45
$ws.Rows.Item(5).RowHeight = 44.45
$ws.Rows.Item(6).RowHeight = 58.35
$ws.Rows.Item(7).RowHeight = 58.35
$ws.Rows.Item(8).RowHeight = 44.45
$ws.Rows.Item(9).RowHeight = 44.45
$ws.Rows.Item(10).RowHeight = 44.45
$ws.Rows.Item(11).RowHeight = 44.45
$ws.Rows.Item(12).RowHeight = 44.45
$ws.Rows.Item(13).RowHeight = 44.45
$ws.Rows.Item(14).RowHeight = 44.45
$ws.Rows.Item(15).RowHeight = 44.45
$ws.Rows.Item(16).RowHeight = 44.45
$ws.Rows.Item(17).RowHeight = 44.45
$ws.Rows.Item(18).RowHeight = 44.45
$ws.Rows.Item(19).RowHeight = 44.45
$ws.Rows.Item(20).RowHeight = 44.45
$ws.Rows.Item(21).RowHeight = 44.45
$ws.Rows.Item(22).RowHeight = 44.45
$ws.Rows.Item(23).RowHeight = 44.45
$ws.Rows.Item(24).RowHeight = 44.45
$ws.Rows.Item(25).RowHeight = 44.45
$ws.Rows.Item(26).RowHeight = 44.45
$ws.Rows.Item(27).RowHeight = 44.45
$ws.Rows.Item(28).RowHeight = 44.45
$ws.Rows.Item(29).RowHeight = 44.45
$ws.Rows.Item(30).RowHeight = 44.45
$ws.Rows.Item(31).RowHeight = 44.45
$ws.Rows.Item(32).RowHeight = 44.45
$ws.Rows.Item(33).RowHeight = 44.45
$ws.Rows.Item(34).RowHeight = 44.45
$ws.Rows.Item(35).RowHeight = 44.45
$ws.Rows.Item(36).RowHeight = 30.4
$ws.Rows.Item(37).RowHeight = 44.45
$ws.Rows.Item(38).RowHeight = 58.35
$ws.Rows.Item(39).RowHeight = 30.4
$ws.Rows.Item(40).RowHeight = 44.45
